# Auto-generated script applying the cryptos.xlsx price/volume update
# (GitHub Actions scheduled refresh of coinranking.com data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.742.57"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "2.269.02"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.643"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "77.10"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.03%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.646"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0969"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("E12").Value = "  -2.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "2.607.05"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.866"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("D17").Value = "2.270.88"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "42.603.51"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "0.0₃0987"
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "234.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.69%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("E27").Value = "  -3.55%  "
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0853"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.98%  "
$ws.Range("E33").Value = "  -4.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.65%  "
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0304"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.207"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "109.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.97%  "
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("B47").Value = "BinanceUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.93%  "
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.12%  "
